$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7647
$ws1.Range("F4").Value = 29
$ws1.Range("F5").Value = 465
$ws1.Range("F6").Value = 4312
$ws1.Range("F7").Value = 328
$ws1.Range("F8").Value = 603
$ws1.Range("F9").Value = 279

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 11

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7647
$ws4.Range("F5").Value = 29
$ws4.Range("F6").Value = 465
$ws4.Range("F7").Value = 4312
$ws4.Range("F8").Value = 328
$ws4.Range("F9").Value = 603
$ws4.Range("F10").Value = 279
$ws4.Range("F12").Value = 11
